$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying analysis changed: the "MuSCs" target cluster (previously one of 4
# per-sending-cluster target rows) was dropped from the receptor-expressing universe,
# so every remaining row is recomputed against a 3-cluster universe and the table
# shrinks from 17 rows (16 data rows) to 13 rows (12 data rows).
$ws.Range("A14:T17").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Vcan"
$ws.Range("C2").Value = "Tlr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = [double]"3"
$ws.Range("F2").Value = [double]"1"
$ws.Range("G2").Value = [double]"3.262296333333333"
$ws.Range("H2").Value = [double]"9.786889"
$ws.Range("I2").Value = [double]"0.01915820289899999"
$ws.Range("J2").Value = [double]"0.01915820289899999"
$ws.Range("K2").Value = [double]"2"
$ws.Range("L2").Value = [double]"0.6666666666666666"
$ws.Range("M2").Value = [double]"1.597802666666666"
$ws.Range("N2").Value = [double]"4.793407999999999"
$ws.Range("O2").Value = [double]"0.02304920886321625"
$ws.Range("P2").Value = [double]"0.02304920886321625"
$ws.Range("Q2").Value = [double]"5.212505780856889"
$ws.Range("R2").Value = [double]"46.912552027712"
$ws.Range("S2").Value = [double]"0.0004415814200629258"
$ws.Range("T2").Value = [double]"0.0004415814200629259"

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Vcan"
$ws.Range("C3").Value = "Tlr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = [double]"3"
$ws.Range("F3").Value = [double]"1"
$ws.Range("G3").Value = [double]"3.262296333333333"
$ws.Range("H3").Value = [double]"9.786889"
$ws.Range("I3").Value = [double]"0.01915820289899999"
$ws.Range("J3").Value = [double]"0.01915820289899999"
$ws.Range("K3").Value = [double]"2"
$ws.Range("L3").Value = [double]"0.6666666666666666"
$ws.Range("M3").Value = [double]"1.27306"
$ws.Range("N3").Value = [double]"3.81918"
$ws.Range("O3").Value = [double]"0.0183646118807784"
$ws.Range("P3").Value = [double]"0.0183646118807784"
$ws.Range("Q3").Value = [double]"4.153098970113333"
$ws.Range("R3").Value = [double]"37.37789073102"
$ws.Range("S3").Value = [double]"0.0003518329605733384"
$ws.Range("T3").Value = [double]"0.0003518329605733385"

# Row 4: ECs -> Resolving-Mac
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Vcan"
$ws.Range("C4").Value = "Tlr2"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = [double]"3"
$ws.Range("F4").Value = [double]"1"
$ws.Range("G4").Value = [double]"3.262296333333333"
$ws.Range("H4").Value = [double]"9.786889"
$ws.Range("I4").Value = [double]"0.01915820289899999"
$ws.Range("J4").Value = [double]"0.01915820289899999"
$ws.Range("K4").Value = [double]"3"
$ws.Range("L4").Value = [double]"1"
$ws.Range("M4").Value = [double]"66.45050433333334"
$ws.Range("N4").Value = [double]"199.351513"
$ws.Range("O4").Value = [double]"0.9585861792560053"
$ws.Range("P4").Value = [double]"0.9585861792560054"
$ws.Range("Q4").Value = [double]"216.7812366347841"
$ws.Range("R4").Value = [double]"1951.031129713057"
$ws.Range("S4").Value = [double]"0.01836478851836373"
$ws.Range("T4").Value = [double]"0.01836478851836373"

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Vcan"
$ws.Range("C5").Value = "Tlr2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = [double]"3"
$ws.Range("F5").Value = [double]"1"
$ws.Range("G5").Value = [double]"145.2141163333334"
$ws.Range("H5").Value = [double]"435.6423490000001"
$ws.Range("I5").Value = [double]"0.8527862647199704"
$ws.Range("J5").Value = [double]"0.8527862647199704"
$ws.Range("K5").Value = [double]"2"
$ws.Range("L5").Value = [double]"0.6666666666666666"
$ws.Range("M5").Value = [double]"1.597802666666666"
$ws.Range("N5").Value = [double]"4.793407999999999"
$ws.Range("O5").Value = [double]"0.02304920886321625"
$ws.Range("P5").Value = [double]"0.02304920886321625"
$ws.Range("Q5").Value = [double]"232.0235023150435"
$ws.Range("R5").Value = [double]"2088.211520835392"
$ws.Range("S5").Value = [double]"0.01965604873121262"
$ws.Range("T5").Value = [double]"0.01965604873121262"

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Vcan"
$ws.Range("C6").Value = "Tlr2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = [double]"3"
$ws.Range("F6").Value = [double]"1"
$ws.Range("G6").Value = [double]"145.2141163333334"
$ws.Range("H6").Value = [double]"435.6423490000001"
$ws.Range("I6").Value = [double]"0.8527862647199704"
$ws.Range("J6").Value = [double]"0.8527862647199704"
$ws.Range("K6").Value = [double]"2"
$ws.Range("L6").Value = [double]"0.6666666666666666"
$ws.Range("M6").Value = [double]"1.27306"
$ws.Range("N6").Value = [double]"3.81918"
$ws.Range("O6").Value = [double]"0.0183646118807784"
$ws.Range("P6").Value = [double]"0.0183646118807784"
$ws.Range("Q6").Value = [double]"184.8662829393133"
$ws.Range("R6").Value = [double]"1663.79654645382"
$ws.Range("S6").Value = [double]"0.015661088768841"
$ws.Range("T6").Value = [double]"0.015661088768841"

# Row 7: FAPs -> Resolving-Mac
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Vcan"
$ws.Range("C7").Value = "Tlr2"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = [double]"3"
$ws.Range("F7").Value = [double]"1"
$ws.Range("G7").Value = [double]"145.2141163333334"
$ws.Range("H7").Value = [double]"435.6423490000001"
$ws.Range("I7").Value = [double]"0.8527862647199704"
$ws.Range("J7").Value = [double]"0.8527862647199704"
$ws.Range("K7").Value = [double]"3"
$ws.Range("L7").Value = [double]"1"
$ws.Range("M7").Value = [double]"66.45050433333334"
$ws.Range("N7").Value = [double]"199.351513"
$ws.Range("O7").Value = [double]"0.9585861792560053"
$ws.Range("P7").Value = [double]"0.9585861792560054"
$ws.Range("Q7").Value = [double]"9649.551266669339"
$ws.Range("R7").Value = [double]"86845.96140002405"
$ws.Range("S7").Value = [double]"0.8174691272199168"
$ws.Range("T7").Value = [double]"0.8174691272199168"

# Row 8: MuSCs -> ECs
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Vcan"
$ws.Range("C8").Value = "Tlr2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = [double]"3"
$ws.Range("F8").Value = [double]"1"
$ws.Range("G8").Value = [double]"21.305189"
$ws.Range("H8").Value = [double]"63.915567"
$ws.Range("I8").Value = [double]"0.1251171236325075"
$ws.Range("J8").Value = [double]"0.1251171236325075"
$ws.Range("K8").Value = [double]"2"
$ws.Range("L8").Value = [double]"0.6666666666666666"
$ws.Range("M8").Value = [double]"1.597802666666666"
$ws.Range("N8").Value = [double]"4.793407999999999"
$ws.Range("O8").Value = [double]"0.02304920886321625"
$ws.Range("P8").Value = [double]"0.02304920886321625"
$ws.Range("Q8").Value = [double]"34.04148779803733"
$ws.Range("R8").Value = [double]"306.3733901823359"
$ws.Range("S8").Value = [double]"0.002883850714970516"
$ws.Range("T8").Value = [double]"0.002883850714970516"

# Row 9: MuSCs -> FAPs
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Vcan"
$ws.Range("C9").Value = "Tlr2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = [double]"3"
$ws.Range("F9").Value = [double]"1"
$ws.Range("G9").Value = [double]"21.305189"
$ws.Range("H9").Value = [double]"63.915567"
$ws.Range("I9").Value = [double]"0.1251171236325075"
$ws.Range("J9").Value = [double]"0.1251171236325075"
$ws.Range("K9").Value = [double]"2"
$ws.Range("L9").Value = [double]"0.6666666666666666"
$ws.Range("M9").Value = [double]"1.27306"
$ws.Range("N9").Value = [double]"3.81918"
$ws.Range("O9").Value = [double]"0.0183646118807784"
$ws.Range("P9").Value = [double]"0.0183646118807784"
$ws.Range("Q9").Value = [double]"27.12278390833999"
$ws.Range("R9").Value = [double]"244.10505517506"
$ws.Range("S9").Value = [double]"0.002297727415150368"
$ws.Range("T9").Value = [double]"0.002297727415150368"

# Row 10: MuSCs -> Resolving-Mac
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Vcan"
$ws.Range("C10").Value = "Tlr2"
$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("E10").Value = [double]"3"
$ws.Range("F10").Value = [double]"1"
$ws.Range("G10").Value = [double]"21.305189"
$ws.Range("H10").Value = [double]"63.915567"
$ws.Range("I10").Value = [double]"0.1251171236325075"
$ws.Range("J10").Value = [double]"0.1251171236325075"
$ws.Range("K10").Value = [double]"3"
$ws.Range("L10").Value = [double]"1"
$ws.Range("M10").Value = [double]"66.45050433333334"
$ws.Range("N10").Value = [double]"199.351513"
$ws.Range("O10").Value = [double]"0.9585861792560053"
$ws.Range("P10").Value = [double]"0.9585861792560054"
$ws.Range("Q10").Value = [double]"1415.740553966986"
$ws.Range("R10").Value = [double]"12741.66498570287"
$ws.Range("S10").Value = [double]"0.1199355455023867"
$ws.Range("T10").Value = [double]"0.1199355455023867"

# Row 11: Resolving-Mac -> ECs
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Vcan"
$ws.Range("C11").Value = "Tlr2"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = [double]"3"
$ws.Range("F11").Value = [double]"1"
$ws.Range("G11").Value = [double]"0.500358"
$ws.Range("H11").Value = [double]"1.501074"
$ws.Range("I11").Value = [double]"0.002938408748521978"
$ws.Range("J11").Value = [double]"0.002938408748521978"
$ws.Range("K11").Value = [double]"2"
$ws.Range("L11").Value = [double]"0.6666666666666666"
$ws.Range("M11").Value = [double]"1.597802666666666"
$ws.Range("N11").Value = [double]"4.793407999999999"
$ws.Range("O11").Value = [double]"0.02304920886321625"
$ws.Range("P11").Value = [double]"0.02304920886321625"
$ws.Range("Q11").Value = [double]"0.7994733466879999"
$ws.Range("R11").Value = [double]"7.195260120192"
$ws.Range("S11").Value = [double]"6.772799697018493E-05"
$ws.Range("T11").Value = [double]"6.772799697018494E-05"

# Row 12: Resolving-Mac -> FAPs
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Vcan"
$ws.Range("C12").Value = "Tlr2"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = [double]"3"
$ws.Range("F12").Value = [double]"1"
$ws.Range("G12").Value = [double]"0.500358"
$ws.Range("H12").Value = [double]"1.501074"
$ws.Range("I12").Value = [double]"0.002938408748521978"
$ws.Range("J12").Value = [double]"0.002938408748521978"
$ws.Range("K12").Value = [double]"2"
$ws.Range("L12").Value = [double]"0.6666666666666666"
$ws.Range("M12").Value = [double]"1.27306"
$ws.Range("N12").Value = [double]"3.81918"
$ws.Range("O12").Value = [double]"0.0183646118807784"
$ws.Range("P12").Value = [double]"0.0183646118807784"
$ws.Range("Q12").Value = [double]"0.6369857554799999"
$ws.Range("R12").Value = [double]"5.73287179932"
$ws.Range("S12").Value = [double]"5.396273621368991E-05"
$ws.Range("T12").Value = [double]"5.396273621368992E-05"

# Row 13: Resolving-Mac -> Resolving-Mac
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Vcan"
$ws.Range("C13").Value = "Tlr2"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = [double]"3"
$ws.Range("F13").Value = [double]"1"
$ws.Range("G13").Value = [double]"0.500358"
$ws.Range("H13").Value = [double]"1.501074"
$ws.Range("I13").Value = [double]"0.002938408748521978"
$ws.Range("J13").Value = [double]"0.002938408748521978"
$ws.Range("K13").Value = [double]"3"
$ws.Range("L13").Value = [double]"1"
$ws.Range("M13").Value = [double]"66.45050433333334"
$ws.Range("N13").Value = [double]"199.351513"
$ws.Range("O13").Value = [double]"0.9585861792560053"
$ws.Range("P13").Value = [double]"0.9585861792560054"
$ws.Range("Q13").Value = [double]"33.249041447218"
$ws.Range("R13").Value = [double]"299.241373024962"
$ws.Range("S13").Value = [double]"0.002816718015338103"
$ws.Range("T13").Value = [double]"0.002816718015338103"

